$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.115.44'
$ws.Range('E2').Value = '  -1.62%  '

# Row 3
$ws.Range('D3').Value = '1.835.75'
$ws.Range('E3').Value = '  -0.76%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.20%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.14%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4623'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.87%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3862'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.16%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07844'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.84%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9608'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.90%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.29%  '

# Row 12
$ws.Range('D12').Value = '1.824.58'
$ws.Range('E12').Value = '  +0.12%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.671'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.59%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.888'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.35%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06855'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.26%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.78%  '

# Row 17
$ws.Range('E17').Value = '  +0.06%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009931'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.88%  '

# Row 19
$ws.Range('E19').Value = '  -2.34%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.01%  '

# Row 21
$ws.Range('D21').Value = '28.130.71'
$ws.Range('E21').Value = '  -1.58%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.293'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.98%  '

# Row 23
$ws.Range('E23').Value = '  -2.23%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.097'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.56%  '

# Row 25
$ws.Range('D25').Value = '2.077.03'
$ws.Range('E25').Value = '  +0.11%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.02%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.73%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.698'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.67%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.973'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.56%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.32%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9354'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.74%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09258'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.81%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.265'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.84%  '

# Row 34
$ws.Range('E34').Value = '  -2.02%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.329'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.26%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05816'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.60%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02116'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.138'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.83%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.708'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.41%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5592'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.08%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.926'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.20%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1753'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.16%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07349'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.50%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.02%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5262'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.16%  '

# Row 46
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.142'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.66%  '

# Row 47
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.136'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.84%  '

# Row 48
$ws.Range('E48').Value = '  -3.61%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '113.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.14%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.002'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.04%  '

# Row 51
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.319'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
